$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1280.3572
$ws.Range("I43").Value = 1045.238
$ws.Range("J43").Value = 1985.7142
$ws.Range("K43").Value = 1045.238
$ws.Range("L43").Value = 1985.7142
$ws.Range("M43").Value = -976.2380000000001
$ws.Range("N43").Value = -2123.7142
$ws.Range("H101").Value = 9291948
$ws.Range("J101").Value = 11364936
$ws.Range("L101").Value = 34094808
$ws.Range("N101").Value = -34098052
$ws.Range("H107").Value = 526.6
$ws.Range("I107").Value = 519.9286
$ws.Range("J107").Value = 620
$ws.Range("K107").Value = 519.9286
$ws.Range("L107").Value = 620
$ws.Range("M107").Value = 1400.0714
$ws.Range("N107").Value = -4460
$ws.Range("H129").Value = 930.1039
$ws.Range("J129").Value = 958.3939
$ws.Range("L129").Value = 2875.1817
$ws.Range("N129").Value = -12875.1817
$ws.Range("H132").Value = 39481788
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 44500
$ws.Range("J133").Value = 44500
$ws.Range("L133").Value = 44500
$ws.Range("N133").Value = -54620
$ws.Range("H135").Value = 9616515
$ws.Range("I135").Value = 9616515
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 86548635
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -86546100
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 2123.8057
$ws.Range("I137").Value = 2056.8572
$ws.Range("J137").Value = 2266.4348
$ws.Range("K137").Value = 6170.571599999999
$ws.Range("L137").Value = 6799.3044
$ws.Range("M137").Value = -3620.571599999999
$ws.Range("N137").Value = -11899.3044
$ws.Range("H138").Value = 4071.7693
$ws.Range("I138").Value = 1617.5625
$ws.Range("J138").Value = 6451.606
$ws.Range("K138").Value = 4852.6875
$ws.Range("L138").Value = 19354.818
$ws.Range("M138").Value = 287.3125
$ws.Range("N138").Value = -29634.818
$ws.Range("H141").Value = 1618.375
$ws.Range("I141").Value = 1192.9333
$ws.Range("K141").Value = 3578.7999
$ws.Range("M141").Value = 1601.2001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 25001228
$ws.Range("I61").Value = 41667540
$ws.Range("J61").Value = 1762.2222
$ws.Range("K61").Value = 41667540
$ws.Range("L61").Value = 1762.2222
$ws.Range("M61").Value = -41667328
$ws.Range("N61").Value = -2186.2222
$ws.Range("H136").Value = 25001228
$ws.Range("I136").Value = 41667540
$ws.Range("J136").Value = 1762.2222
$ws.Range("K136").Value = 125002620
$ws.Range("L136").Value = 5286.6666
$ws.Range("M136").Value = -125000070
$ws.Range("N136").Value = -10386.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 901372.6
$ws.Range("I22").Value = 1228900.9
$ws.Range("K22").Value = 1228900.9
$ws.Range("M22").Value = -1228727.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 34483780
$ws.Range("I58").Value = 40000748
$ws.Range("J58").Value = 2725
$ws.Range("K58").Value = 40000748
$ws.Range("L58").Value = 2725
$ws.Range("M58").Value = -40000545
$ws.Range("N58").Value = -3131
$ws.Range("H94").Value = 6102
$ws.Range("I94").Value = 28175
$ws.Range("J94").Value = 1455.0526
$ws.Range("K94").Value = 28175
$ws.Range("L94").Value = 1455.0526
$ws.Range("M94").Value = -27724
$ws.Range("N94").Value = -2357.0526
$ws.Range("H121").Value = 120000
$ws.Range("J121").Value = 120000
$ws.Range("L121").Value = 120000
$ws.Range("N121").Value = -122620
$ws.Range("H134").Value = 9434826
$ws.Range("I134").Value = 670.1957
$ws.Range("K134").Value = 2010.5871
$ws.Range("M134").Value = 524.4129
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 34483780
$ws.Range("I136").Value = 40000748
$ws.Range("J136").Value = 2725
$ws.Range("K136").Value = 120002244
$ws.Range("L136").Value = 8175
$ws.Range("M136").Value = -119999694
$ws.Range("N136").Value = -13275
$ws.Range("H138").Value = 50780
$ws.Range("J138").Value = 50780
$ws.Range("L138").Value = 50780
$ws.Range("N138").Value = -61060
$ws.Range("H140").Value = 40580
$ws.Range("J140").Value = 40580
$ws.Range("L140").Value = 40580
$ws.Range("N140").Value = -50940
$ws.Range("H141").Value = 29171.545
$ws.Range("J141").Value = 29260.857
$ws.Range("L141").Value = 29260.857
$ws.Range("N141").Value = -39620.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 6509
$ws.Range("I124").Value = 2476.6667
$ws.Range("J124").Value = 9533.25
$ws.Range("K124").Value = 7430.000100000001
$ws.Range("L124").Value = 28599.75
$ws.Range("M124").Value = -2520.000100000001
$ws.Range("N124").Value = -38419.75
$ws.Range("H131").Value = 923.05
$ws.Range("J131").Value = 951.16486
$ws.Range("L131").Value = 2853.49458
$ws.Range("N131").Value = -12933.49458
$ws.Range("H137").Value = 26318392
$ws.Range("I137").Value = 100001300
$ws.Range("J137").Value = 3067.5715
$ws.Range("K137").Value = 300003900
$ws.Range("L137").Value = 9202.7145
$ws.Range("M137").Value = -299998800
$ws.Range("N137").Value = -19402.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1856.5
$ws.Range("I7").Value = 1475.3334
$ws.Range("J7").Value = 3000
$ws.Range("K7").Value = 1475.3334
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = -1363.3334
$ws.Range("N7").Value = -3224
$ws.Range("H46").Value = 1487.3572
$ws.Range("I46").Value = 832.1667
$ws.Range("J46").Value = 1978.75
$ws.Range("K46").Value = 832.1667
$ws.Range("L46").Value = 1978.75
$ws.Range("M46").Value = -644.1667
$ws.Range("N46").Value = -2354.75
$ws.Range("H126").Value = 1856.5
$ws.Range("I126").Value = 1475.3334
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4426.0002
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1956.0002
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 25979400
$ws.Range("I132").Value = 40818612
$ws.Range("J132").Value = 10781.75
$ws.Range("K132").Value = 122455836
$ws.Range("L132").Value = 32345.25
$ws.Range("M132").Value = -122453306
$ws.Range("N132").Value = -37405.25
$ws.Range("H136").Value = 68163230
$ws.Range("I136").Value = 77383140
$ws.Range("J136").Value = 52634956
$ws.Range("K136").Value = 232149420
$ws.Range("L136").Value = 157904868
$ws.Range("M136").Value = -232146870
$ws.Range("N136").Value = -157909968

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2123.125
$ws.Range("I126").Value = 2046.9445
$ws.Range("J126").Value = 2351.6667
$ws.Range("K126").Value = 6140.833500000001
$ws.Range("L126").Value = 7055.000100000001
$ws.Range("M126").Value = -3670.833500000001
$ws.Range("N126").Value = -11995.0001
$ws.Range("H132").Value = 25943.17
$ws.Range("I132").Value = 34747.582
$ws.Range("J132").Value = 8884.625
$ws.Range("K132").Value = 104242.746
$ws.Range("L132").Value = 26653.875
$ws.Range("M132").Value = -101712.746
$ws.Range("N132").Value = -31713.875
$ws.Range("H138").Value = 47750
$ws.Range("J138").Value = 47750
$ws.Range("L138").Value = 47750
$ws.Range("N138").Value = -58030

